# Apply updated odds values to Sheet1, matching the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("K2").Value = 1.87
$ws.Range("L2").Value = 3.5
$ws.Range("U2").Value = 2.1
$ws.Range("V2").Value = 1.67
$ws.Range("W2").Value = 7
$ws.Range("Z2").Value = 29
$ws.Range("AK2").Value = 29
$ws.Range("AL2").Value = 26

# Row 5
$ws.Range("L5").Value = 1.92
$ws.Range("Q5").Value = 1.62
$ws.Range("U5").Value = 1.8
$ws.Range("V5").Value = 1.8

# Row 7
$ws.Range("G7").Value = 3.1
$ws.Range("I7").Value = 2.2
$ws.Range("M7").Value = 1.03
$ws.Range("O7").Value = 1.25
$ws.Range("U7").Value = 1.72
$ws.Range("X7").Value = 15
$ws.Range("AA7").Value = 23
$ws.Range("AB7").Value = 29
$ws.Range("AE7").Value = 13
$ws.Range("AH7").Value = 8.5
$ws.Range("AU7").Value = 7.5

# Row 9
$ws.Range("G9").Value = 1.48
$ws.Range("I9").Value = 5.75
$ws.Range("J9").Value = 1.95
$ws.Range("AD9").Value = 9
$ws.Range("AR9").Value = 41
$ws.Range("AW9").Value = 7.5
$ws.Range("AZ9").Value = 81

# Row 10
$ws.Range("M10").Value = 1.08
$ws.Range("N10").Value = 8
